# Updated cryptos list with GitHub Actions — refreshed Price/Volume(1h)
# columns, plus the bottom-of-list coin rotation (RocketPoolETH dropped,
# TheGraph/BEAM shift up, FraxShare newly added as row 51).
#
# Price values in column D are plain numeric-looking text (e.g. "354.33")
# that must stay literal strings, not be auto-coerced to numbers by the
# Value setter, so those cells are pre-formatted as Text ("@") first.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.855.88"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.978.96"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.33"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.33"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.66"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0850"
$ws.Range("E12").Value = "  -3.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.92"
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("D14").Value = "3.459.58"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.51"
$ws.Range("D16").Value = "2.966.98"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.992"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "51.812.36"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.35"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.76"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.52"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("E25").Value = "  -4.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.176"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.70"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.108"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.23"
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.09"
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.62"
$ws.Range("E33").Value = "  -8.13%  "
$ws.Range("E34").Value = "  +12.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.90"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0428"
$ws.Range("E36").Value = "  -3.72%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.33"
$ws.Range("E41").Value = "  -6.50%  "
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.94"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.74"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "2.111.42"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.27"
$ws.Range("E47").Value = "  -5.33%  "
$ws.Range("E48").Value = "  -7.35%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.239"
$ws.Range("E49").Value = "  -4.08%  "
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0331"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.91"
$ws.Range("E51").Value = "  -3.64%  "
